# The workbook currently has three sheets:
#   strategy_id-0, strategy_id-5004, strategy_id-5008
# The target state renames "strategy_id-5008" to "strategy_id-5007" and adds
# a new sheet "strategy_id-5009" (an exact duplicate of that sheet's data)
# right after it.

$wb = $excel.ActiveWorkbook

# 1. Rename strategy_id-5008 -> strategy_id-5007
$ws = $wb.Worksheets.Item("strategy_id-5008")
$ws.Name = "strategy_id-5007"

# 2. Duplicate it (placed immediately after) and rename the copy to
#    strategy_id-5009
$ws.Copy([System.Reflection.Missing]::Value, $ws)
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "strategy_id-5009"

# Copying activates the new sheet; restore the original active sheet so we
# don't introduce an unrelated view-state change.
$wb.Worksheets.Item("strategy_id-0").Activate()
